$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.738.12'
$ws.Range("E2").Value = '  -1.16%  '

$ws.Range("D3").Value = '3.405.06'
$ws.Range("E3").Value = '  -0.99%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.46'
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.48'
$ws.Range("E6").Value = '  -4.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +7.08%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.732'
$ws.Range("E9").Value = '  +6.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  +10.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.58'
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.141'
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.00'
$ws.Range("E13").Value = '  +5.91%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.949.84'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.04'
$ws.Range("E15").Value = '  +5.22%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000206'
$ws.Range("E16").Value = '  +49.60%  '

$ws.Range("D17").Value = '3.382.23'
$ws.Range("E17").Value = '  -2.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.06'
$ws.Range("E18").Value = '  +5.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.07'
$ws.Range("E19").Value = '  +4.32%  '

$ws.Range("D20").Value = '61.756.41'
$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.71'
$ws.Range("E21").Value = '  +41.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.13'
$ws.Range("E22").Value = '  +8.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.15'
$ws.Range("E23").Value = '  -1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +2.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.33'
$ws.Range("E26").Value = '  +11.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.71'
$ws.Range("E27").Value = '  +6.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.74'
$ws.Range("E28").Value = '  -0.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.57'
$ws.Range("E29").Value = '  -1.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.73'
$ws.Range("E30").Value = '  -2.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.92'
$ws.Range("E31").Value = '  +4.32%  '

$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.167'
$ws.Range("E32").Value = '  -3.87%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.113'
$ws.Range("E33").Value = '  -2.63%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.51'
$ws.Range("E34").Value = '  -1.37%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0495'
$ws.Range("E36").Value = '  +1.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.49'
$ws.Range("E37").Value = '  +3.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.36'
$ws.Range("E39").Value = '  -1.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.134'
$ws.Range("E40").Value = '  +6.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.89'
$ws.Range("E41").Value = '  -1.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.312'
$ws.Range("E42").Value = '  -4.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.04'
$ws.Range("E43").Value = '  +1.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.14'
$ws.Range("E44").Value = '  +3.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("E46").Value = '  +7.65%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.54'
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.23'
$ws.Range("E48").Value = '  +3.05%  '

$ws.Range("D49").Value = '3.757.23'
$ws.Range("E49").Value = '  -0.83%  '

$ws.Range("D50").Value = '2.102.97'
$ws.Range("E50").Value = '  -1.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '105.03'
$ws.Range("E51").Value = '  +25.55%  '
